# Generate Report for Handback
# This script regenerates the localization-status report for a handback:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - The per-language sheets (zh-cn, de-de) gain "Latest Target File" (E) and
#    "Latest Handback File" (F) hyperlink entries pointing at the handed-back files
#  - "Latest Handback DateTime" (G) is stamped with the handback timestamp

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

# --- Overview sheet: refresh status text so the workbook stays internally
# consistent (it shares the same status string as the per-language sheets). ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusNew
$overview.Range("C2").Value = $statusNew
$overview.Range("B3").Value = $statusNew
$overview.Range("C3").Value = $statusNew

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B2").Value = $statusNew
$zh.Range("B3").Value = $statusNew

# Latest Target File / Latest Handback File hyperlinks for rows 2 and 3.
# Both rows reference the same handed-back md/xlf pair (the file that was
# actually processed in this handback batch).
$zhMdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/591f9f87cf93ecc51acfc20738d14a0a81655c5f/e2e/a28591f3-d3b4-45d5-86e2-cebe1e59fd36.md"
$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/969fb623a3612c61acf3c447c540bb0a8d227084/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.zh-cn.xlf"

$zh.Hyperlinks.Add($zh.Range("E2"), $zhMdAddress, "", "", "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F2"), $zhXlfAddress, "", "", "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E3"), $zhMdAddress, "", "", "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F3"), $zhXlfAddress, "", "", "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.zh-cn.xlf") | Out-Null

$zh.Range("G2").Value = "2016-03-10 19:06:40"
$zh.Range("G3").Value = "2016-03-10 19:06:40"

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")

$de.Range("B2").Value = $statusNew
$de.Range("B3").Value = $statusNew

$deMdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/591f9f87cf93ecc51acfc20738d14a0a81655c5f/e2e/a28591f3-d3b4-45d5-86e2-cebe1e59fd36.md"
$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a8a93091275577b9bfbd16594fef9af79cc9368d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.de-de.xlf"

$de.Hyperlinks.Add($de.Range("E2"), $deMdAddress, "", "", "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F2"), $deXlfAddress, "", "", "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E3"), $deMdAddress, "", "", "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F3"), $deXlfAddress, "", "", "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.de-de.xlf") | Out-Null

$de.Range("G2").Value = "2016-03-10 19:06:51"
$de.Range("G3").Value = "2016-03-10 19:06:51"
